$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data rows (row 2 and row 3), columns B..J. Column A (trial index) stays put.
for ($col = 2; $col -le 10; $col++) {
    $cell2 = $ws.Cells.Item(2, $col)
    $cell3 = $ws.Cells.Item(3, $col)
    $v2 = $cell2.Value2
    $v3 = $cell3.Value2
    $cell2.Value = $v3
    $cell3.Value = $v2
}

# Update the active selection to match the workbook's saved view state.
$ws.Range("B2:K3").Select()
